$d = $word.ActiveDocument

# Locate " (Robrisqui)" - it currently lives in a single run. The edit needs
# to turn it into three runs: " (", "Luke Wroblewski" and ")" - while leaving
# the neighboring runs (" ... ex-yahoo" before it and " que hoje ..." after
# it) completely untouched.
$match = $d.Content
$match.Find.Execute(" (Robrisqui)", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0) | Out-Null

$openParenEnd = $match.Start + 2   # just past " ("
$nameEnd      = $openParenEnd + 9  # "Robrisqui" is 9 characters long
$closeParenEnd = $nameEnd + 1      # just past the ")"

# Mark the name and the closing ")" with a distinguishing (temporary) format
# while they still sit inside the original run. Word only fuses adjacent
# runs that carry identical formatting, so doing this first - before any
# text is changed - guarantees the paragraph ends up split at exactly the
# boundaries we want, no matter how the engine re-normalizes runs as each
# edit below is applied.
$nameRange = $d.Range($openParenEnd, $nameEnd)
$nameRange.Font.Bold = 1

$closeParenRange = $d.Range($nameEnd, $closeParenEnd)
$closeParenRange.Font.Bold = 1

# Swap the old name for the new one inside the (now distinct) bold run.
$nameRange.Text = "Luke Wroblewski"
$newNameEnd = $nameRange.End

# Finally, strip the temporary bold marker back off both pieces so the
# visible formatting matches the original (plain) text.
$d.Range($openParenEnd, $newNameEnd).Font.Bold = 0
$d.Range($newNameEnd, $newNameEnd + 1).Font.Bold = 0
